$d = $word.ActiveDocument

# --- Edit 1: restructure the phishing-detection paragraph (was paragraph 4) ---
$p4 = $d.Paragraphs(4)
$r4 = $p4.Range.Duplicate
$xmlPara4 = @'
<w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">For automatic detection of </w:t></w:r><w:r><w:t xml:space="preserve">suspicious </w:t></w:r><w:r><w:t xml:space="preserve">emails there is a number o approaches to detect phishing emails using email headers, SMTP headers NIDS logs and others. What we are interested in is the approaches that are focused on processing email content or parts of it. One way of doing that is preforming authorship of identification, to analyze email features and assuring that who claimed to be </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>the sender of the email in the content is really the sender.</w:t></w:r><w:r><w:t xml:space="preserve"> These approaches characterize emails by using statistics on frequency of words using n-grams. One example of that is ASCAI which generates a writeprint of senders, then generate another writeprint for the current email sender, and compare against the list of known senders’ writeprints to verify authorship.</w:t></w:r></w:p>
'@
$r4.InsertXML($xmlPara4)

# --- Edit 2: insert the new Introduction section before 'Related work' (paragraph 1) ---
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$r1.InsertParagraphBefore()
$newp = $d.Paragraphs(1)
$newr = $newp.Range
$xmlIntro = @'
<w:p><w:pPr><w:jc w:val="center"/></w:pPr><w:r><w:t>Introduction</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">One of the most dangerous threats to the information system at all is social engineering. The emotional manipulation for the target. The weakest part of an information system often is the human factor. In most cases for the social engineer hacking the emotional state of people is much easier that the hardened computer systems. </w:t></w:r><w:r><w:t xml:space="preserve">Any type of data that the social engineer can collect from a normal conversation with a human victim can be relevant for him </w:t></w:r><w:r><w:t>to gain access to his target like guessing some personal password for example.</w:t></w:r><w:r><w:t xml:space="preserve"> The use of modern technological devices has been of a great benefit for the attacker to gain more access and reach for his target.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Social engineering attacks depends on a some type of communication between the attacker and a worker for example in order for the attack to collect sensitive data </w:t></w:r><w:r><w:t>filter</w:t></w:r><w:r><w:t>ed</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">from the conversation, that data can be of a great help for the attacker to guess a password, a credit-card number ,or know a name of a coworker which can identify that someone is a potential target for the attacker to operate on with chain attacks</w:t></w:r><w:r><w:t xml:space="preserve"> and in both cases this also might support larger attacks</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> For the attacker to get his target’s data, he would convince the victim to do some action like going to a suspicious website or a fake website to enter his confidential details. Other types of that would require the attacker to impersonate a character that the user can trust and pretend it is him on an SMS, Email, or even a phone call, and the victim will be easily manipulated.</w:t></w:r><w:r><w:t xml:space="preserve"> Surprisingly this is effective and encourage the attacker to use it more frequently.</w:t></w:r></w:p><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Phishing attacks are a type of social engineering attacks that depend on </w:t></w:r><w:r><w:t>disguising</w:t></w:r><w:r><w:t xml:space="preserve"> as a trustworthy entity on the internet and deceive people to gain their sensitive information</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/></w:pPr></w:p>
'@
$newr.InsertXML($xmlIntro)

Write-Host "Edit complete"
